$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add a new worksheet at the very end of the workbook (becomes
#    "Sheet25" automatically, following the existing Sheet1..Sheet24
#    naming). It will become the active/selected tab, mirroring what
#    Sheet19 used to be.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)

# ---------------------------------------------------------------------
# 2. Populate the new sheet with the same shape/content as Sheet20,
#    except for A5 which references a new "eurostat" source instead of
#    "umar", and rows 2-5 only going out to column I (no table_code /
#    series_code values for this data set).
# ---------------------------------------------------------------------
$newSheet.Cells.Item(1, 1).Value = "source"
$newSheet.Cells.Item(1, 2).Value = "author"
$newSheet.Cells.Item(1, 3).Value = "table_name"
$newSheet.Cells.Item(1, 4).Value = "dimensions"
$newSheet.Cells.Item(1, 5).Value = "dimension_levels_text"
$newSheet.Cells.Item(1, 6).Value = "dimension_levels_code"
$newSheet.Cells.Item(1, 7).Value = "unit"
$newSheet.Cells.Item(1, 8).Value = "interval"
$newSheet.Cells.Item(1, 9).Value = "series_name"
$newSheet.Cells.Item(1, 10).Value = "table_code"
$newSheet.Cells.Item(1, 11).Value = "series_code"

$newSheet.Cells.Item(2, 1).Value = "umar"
$newSheet.Cells.Item(2, 2).Value = "mz"
$newSheet.Cells.Item(2, 3).Value = "dfg"
$newSheet.Cells.Item(2, 4).Value = "dff"
$newSheet.Cells.Item(2, 5).Value = 3
$newSheet.Cells.Item(2, 6).Value = 234
$newSheet.Cells.Item(2, 7).Value = "%"
$newSheet.Cells.Item(2, 8).Value = "M"
$newSheet.Cells.Item(2, 9).Value = "tesx"

$newSheet.Cells.Item(3, 1).Value = "umar"
$newSheet.Cells.Item(3, 2).Value = "mz"
$newSheet.Cells.Item(3, 3).Value = "dfg"
$newSheet.Cells.Item(3, 4).Value = "dff"
$newSheet.Cells.Item(3, 5).Value = 2
$newSheet.Cells.Item(3, 6).Value = 1123
$newSheet.Cells.Item(3, 7).Value = "%"
$newSheet.Cells.Item(3, 8).Value = "M"
$newSheet.Cells.Item(3, 9).Value = "sdt"

$newSheet.Cells.Item(4, 1).Value = "umar"
$newSheet.Cells.Item(4, 2).Value = "mz"
$newSheet.Cells.Item(4, 3).Value = "kjh"
$newSheet.Cells.Item(4, 4).Value = "dff"
$newSheet.Cells.Item(4, 5).Value = 4
$newSheet.Cells.Item(4, 6).Value = 1
$newSheet.Cells.Item(4, 7).Value = "%"
$newSheet.Cells.Item(4, 8).Value = "M"

$newSheet.Cells.Item(5, 1).Value = "eurostat"
$newSheet.Cells.Item(5, 2).Value = "mz"
$newSheet.Cells.Item(5, 3).Value = "kjh"
$newSheet.Cells.Item(5, 4).Value = "dff"
$newSheet.Cells.Item(5, 5).Value = 5
$newSheet.Cells.Item(5, 6).Value = 12
$newSheet.Cells.Item(5, 7).Value = "%"
$newSheet.Cells.Item(5, 8).Value = "M"

# ---------------------------------------------------------------------
# 3. Make the new sheet the active tab/selection (this also clears
#    tabSelected on whichever sheet used to carry it).
# ---------------------------------------------------------------------
$newSheet.Activate()
$newSheet.Range("A6").Select()

# ---------------------------------------------------------------------
# 4. On "Sheet20", grow the current selection out to the full sheet
#    (Ctrl+A / Select All), instead of just the single cell K17.
# ---------------------------------------------------------------------
$ws20 = $wb.Worksheets.Item("Sheet20")
$ws20.Activate()
$ws20.Range("K17").Select()
$ws20.Cells.Select()

# Leave the previously-active sheet (now Sheet25) focused, matching the
# new activeTab in the saved workbook.
$newSheet.Activate()
